$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections in existing rows ---
$ws.Range("C2").Value = "Full Stack Engineer"
$ws.Range("B5").Value = "Dhanush"

# --- New employee row (emp033 / Koushik / Intern) ---
$ws.Range("A34").Value = "emp033"
$ws.Range("B34").Value = "Koushik"
$ws.Range("C34").Value = "Intern"
$ws.Range("D34").Formula = '=_xlfn.CONCAT(A34,"@emp.com")'

# --- New "Action" column (E) with a repeating Add/Update/Delete tag per row ---
$ws.Range("E1").Value = "Action"

$actions = @("Add", "Update", "Delete")
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 5).Value = $actions[($r - 2) % 3]
}

# --- View state: show formulas instead of values, move selection to C5 ---
$ws.Range("C5").Select()
$excel.ActiveWindow.DisplayFormulas = $true
